# Fruta / hortaliza, semanal
# Insert two new weekly rows (Provincia de Curicó, week of 2023-05-16) above
# the existing row 32, pushing the previous rows 32-41 down to 34-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 32 (formatting/style is inherited from the
# row above, so the date column keeps its date number format).
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

# New row 32: Granada, Wonderfull, Especial, Provincia de Curicó
$ws.Range("A32").Value = 9
$ws.Range("B32").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 45062
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100104
$ws.Range("H32").Value = "Frutos de pepita"
$ws.Range("I32").Value = 100104001
$ws.Range("J32").Value = "Granada"
$ws.Range("K32").Value = "Wonderfull"
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 10500
$ws.Range("O32").Value = 10500
$ws.Range("P32").Value = 10500
$ws.Range("Q32").Value = "$/caja 15 kilos granel"
$ws.Range("R32").Value = "Provincia de Curicó"
$ws.Range("S32").Value = 700
$ws.Range("T32").Value = 15

# New row 33: Granada, Wonderfull, Primera, Provincia de Curicó
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 45062
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100104
$ws.Range("H33").Value = "Frutos de pepita"
$ws.Range("I33").Value = 100104001
$ws.Range("J33").Value = "Granada"
$ws.Range("K33").Value = "Wonderfull"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 9000
$ws.Range("O33").Value = 9000
$ws.Range("P33").Value = 9000
$ws.Range("Q33").Value = "$/caja 15 kilos granel"
$ws.Range("R33").Value = "Provincia de Curicó"
$ws.Range("S33").Value = 600
$ws.Range("T33").Value = 15
